# netCrypto.xlsx - "Add files via upload" re-save.
#
# The workbook was re-opened/re-saved (from a different machine/user -
# "DanB"), the viewport was scrolled and a new cell selected, and the
# USD Amount in T2 was corrected from 627406 to 647200.
#
# Reproduce the user-visible / model-level changes through the Excel
# object model:
#   - sheetView.topLeftCell  H1 -> L1   (ActiveWindow scroll position)
#   - selection               T3 -> V11 (Range.Select)
#   - T2 value                627406 -> 647200

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the window so column L / row 1 is the top-left visible cell
# (was H1).
$excel.ActiveWindow.ScrollColumn = 12   # column L
$excel.ActiveWindow.ScrollRow = 1       # row 1

# Move the selection to V11 (was T3).
$ws.Range("V11").Select() | Out-Null

# Correct the USD Amount value in T2.
$ws.Range("T2").Value = 647200
